# Updated: st 30. 06. 2021
# Apply corrected AgTests (col F) / AgPosit (col G) figures for late-June 2021 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F423").Value = 437381
$ws.Range("G423").Value = 635
$ws.Range("F425").Value = 137408
$ws.Range("F426").Value = 106626
$ws.Range("F427").Value = 88933
$ws.Range("F428").Value = 102252
$ws.Range("G428").Value = 386
$ws.Range("F429").Value = 171078
$ws.Range("F430").Value = 169244
$ws.Range("F432").Value = 118108
$ws.Range("G432").Value = 413
$ws.Range("F433").Value = 85888
$ws.Range("F434").Value = 78983
$ws.Range("F435").Value = 83049
$ws.Range("F436").Value = 139115
$ws.Range("F438").Value = 118271
$ws.Range("F439").Value = 86574
$ws.Range("F440").Value = 72802
$ws.Range("F441").Value = 65740
$ws.Range("F442").Value = 67231
$ws.Range("F443").Value = 102667
$ws.Range("F444").Value = 99951
$ws.Range("F446").Value = 86550
$ws.Range("F447").Value = 64575
$ws.Range("F448").Value = 58739
$ws.Range("F449").Value = 59676
$ws.Range("F450").Value = 87448
$ws.Range("G450").Value = 165
$ws.Range("F451").Value = 82308
$ws.Range("F453").Value = 67237
$ws.Range("F454").Value = 50694
$ws.Range("F455").Value = 49957
$ws.Range("F456").Value = 47752
$ws.Range("F457").Value = 75246
$ws.Range("F458").Value = 67376
$ws.Range("F460").Value = 55574
$ws.Range("F461").Value = 43426
$ws.Range("F462").Value = 41893
$ws.Range("F463").Value = 44744
$ws.Range("G463").Value = 67
$ws.Range("F464").Value = 69718
$ws.Range("F465").Value = 58087
$ws.Range("F467").Value = 50088
$ws.Range("F468").Value = 40454
$ws.Range("F469").Value = 39092
$ws.Range("F470").Value = 41413
$ws.Range("F471").Value = 62447
$ws.Range("G471").Value = 50
$ws.Range("F472").Value = 47431
$ws.Range("F473").Value = 38647
$ws.Range("G473").Value = 167
$ws.Range("F474").Value = 43562
$ws.Range("F475").Value = 33698
$ws.Range("G475").Value = 25
$ws.Range("F476").Value = 34680
$ws.Range("F477").Value = 36470
$ws.Range("F478").Value = 49965
$ws.Range("F479").Value = 37297
$ws.Range("G479").Value = 31
$ws.Range("F480").Value = 30392
$ws.Range("F481").Value = 39429
$ws.Range("G481").Value = 24
